$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1) and "全部类型" (index 4) both carry the same
# per-event rows; "全部类型" has one extra row (the single 演出 entry)
# inserted at row 33, which shifts everything below it down by one.

$sheet1 = $wb.Worksheets.Item(1)   # 展览
$sheet4 = $wb.Worksheets.Item(4)   # 全部类型

# row -> new F (想去人数) value, keyed by row number in 展览 (sheet1)
$sheet1Updates = @{
    3  = 341
    4  = 437
    5  = 1743
    6  = 89
    7  = 2192
    11 = 4968
    12 = 12
    15 = 228
    18 = 37
    21 = 3924
    23 = 682
    26 = 108
    28 = 26
    30 = 92
    34 = 986
    35 = 2504
    36 = 428
    37 = 6
}

# row -> new F value, keyed by row number in 全部类型 (sheet4, shifted by +1 past row 33)
$sheet4Updates = @{
    3  = 341
    4  = 437
    5  = 1743
    6  = 89
    7  = 2192
    11 = 4968
    12 = 12
    15 = 228
    18 = 37
    21 = 3924
    23 = 682
    26 = 108
    28 = 26
    30 = 92
    35 = 986
    36 = 2504
    37 = 428
    38 = 6
}

foreach ($row in $sheet1Updates.Keys) {
    $sheet1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

foreach ($row in $sheet4Updates.Keys) {
    $sheet4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}

# Event "宜春·逆光ZERO动漫游戏展" (row 19 in both sheets) was cancelled:
# name gets a "（取消）" suffix and the lowest-price cell becomes "不可售".
$sheet1.Cells.Item(19, 3).Value = "宜春·逆光ZERO动漫游戏展（取消）"
$sheet1.Cells.Item(19, 7).Value = "不可售"

$sheet4.Cells.Item(19, 3).Value = "宜春·逆光ZERO动漫游戏展（取消）"
$sheet4.Cells.Item(19, 7).Value = "不可售"
